$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current column C holds the "Nom" / nomPatCandidat column (header row1 = hideCols label,
# row2 = formula text). We insert a brand-new column D for "Nom usuel", shifting every
# column from the old D onward one position to the right.
$ws.Columns("D").Insert()

# The inserted column should carry the same look & feel (width/style) as column C, the
# column it now sits next to, matching how Excel normally extends the neighbouring format
# when a column is inserted.
$ws.Columns("D").ColumnWidth = $ws.Columns("C").ColumnWidth
$ws.Range("D1").Style = $ws.Range("C1").Style
$ws.Range("D2").Style = $ws.Range("C2").Style
$ws.Range("D3").Style = $ws.Range("C3").Style

# Row 1 (s=2 header style): rename the old "Nom" hideCols header to "Nom patronymique"
# (column C) and add the new "Nom usuel" hideCols header (column D).
$ws.Range("C1").Value = '<jt:hideCols test="${nomPatHide}">Nom patronymique</jt:hideCols>'
$ws.Range("D1").Value = '<jt:hideCols test="${nomUsuHide}">Nom usuel</jt:hideCols>'

# Row 2 (s=3 data/template style): column C keeps its existing nomPatCandidat expression,
# column D gets the new nomUsuCandidat expression.
$ws.Range("D2").Value = '${cand.candidat.nomUsuCandidat}'
